$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$g2Value = @"
<get>
    <filter type="subtree">
      <components xmlns="http://openconfig.net/yang/platform">
        <component>
          <name>OCH-1-1</name>
             <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
               <config>
                 <frequency></frequency>
               </config>
             </optical-channel>
        </component>
      </components>
    </filter>
</get>
"@

$j2Value = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:82560579-297c-4f64-bd8f-b47ad42417b5"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>OCH-1-1</name>
    <optical-channel xmlns="http://openconfig.net/yang/terminal-device">
     <config>
      <frequency>192000000</frequency>
     </config>
    </optical-channel>
   </component>
  </components>
 </data>
</rpc-reply>
"@

$ws.Range("G2").Value = $g2Value.TrimEnd("`r", "`n")
$ws.Range("J2").Value = $j2Value.TrimEnd("`r", "`n")
